# MonsterTable.xlsx edit
# - Rename column I header (row2) from the old "hud y offset (pixel)" label
#   to the new ratio-based label "hud y 좌표(1=>iso tile height)"
# - Rewrite the f_EffectOffset (H) / f_HUDOffset (I) numeric columns from the
#   old fixed pixel values (15 / 100) to new normalized ratio values used for
#   the attack-particle angle calculation / hud raycast offset.
# - Widen column I so the new (longer) header text fits, matching Excel's
#   "best fit" auto-resize behaviour.
# - Update the active selection to I3 (matches the saved view state).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text (row 2) ---
$ws.Range("I2").Value = "hud y 좌표(1=>iso tile height)"

# --- New H/I values per monster row (3-15) ---
$ws.Range("H3").Value = 0.2
$ws.Range("I3").Value = 0.8

$ws.Range("H4").Value = 0.15
$ws.Range("I4").Value = 0.9

$ws.Range("H5").Value = 0.4
$ws.Range("I5").Value = 1

$ws.Range("H6").Value = 0.4
$ws.Range("I6").Value = 1

$ws.Range("H7").Value = 0.4
$ws.Range("I7").Value = 1

$ws.Range("H8").Value = 0.2
$ws.Range("I8").Value = 0.9

$ws.Range("H9").Value = 0.5
$ws.Range("I9").Value = 1.15

$ws.Range("H10").Value = 0.2
$ws.Range("I10").Value = 0.9

$ws.Range("H11").Value = 0.2
$ws.Range("I11").Value = 0.9

$ws.Range("H12").Value = 0.5
$ws.Range("I12").Value = 1.15

$ws.Range("H13").Value = 0.2
$ws.Range("I13").Value = 0.9

$ws.Range("H14").Value = 0.2
$ws.Range("I14").Value = 0.9

$ws.Range("H15").Value = 0.2
$ws.Range("I15").Value = 0.9

# --- Column width (column I widened to fit new header) ---
$ws.Columns.Item(9).ColumnWidth = 29.375

# --- View / selection state ---
$ws.Range("I3").Select()
